# Fruta / hortaliza, semanal
#
# The weekly refresh drops 3 brand-new price rows at the top of the
# "Comercializadora del Agro de Limari - Frutilla" block (row 245) and
# pushes the rest of the block (old rows 245-319) down by 3 rows so the
# whole history shifts (old row 245 -> new row 248, ..., old row 319 ->
# new row 322). The inserted rows share every "constant" column with the
# rest of the block (market/region/product metadata, unit, origin, kg
# factor) and only differ in Fecha (D), Volumen (M), Precio minimo (N),
# Precio maximo (O), Precio promedio ponderado (P) and Precio $/Kg (S).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing data down by inserting 3 blank rows right above the
# first row of the block (row 245). Excel moves all cell content/styles
# currently in rows 245:319 down by 3, to 248:322, and keeps the
# worksheet dimension (A1:T322) in sync automatically.
$ws.Rows("245:247").Insert()

# The row that used to be 245 is now 248 - reuse it as the template for
# every column that stays constant across the whole block.
$constCols = @("A","B","C","E","F","G","H","I","J","K","Q","R","T")
foreach ($col in $constCols) {
    $templateValue = $ws.Range($col + "248").Value2
    $ws.Range($col + "245").Value = $templateValue
    $ws.Range($col + "246").Value = $templateValue
    $ws.Range($col + "247").Value = $templateValue
}

# Calidad cycles Especial/Primera/Segunda every 3 rows.
$ws.Range("L245").Value = "Especial"
$ws.Range("L246").Value = "Primera"
$ws.Range("L247").Value = "Segunda"

# New sample date for the 3 inserted rows.
$ws.Range("D245").Value = 44588
$ws.Range("D246").Value = 44588
$ws.Range("D247").Value = 44588

# Row 245 - Especial
$ws.Range("M245").Value = 400
$ws.Range("N245").Value = 11500
$ws.Range("O245").Value = 12000
$ws.Range("P245").Value = 11750
$ws.Range("S245").Value = 1679

# Row 246 - Primera
$ws.Range("M246").Value = 500
$ws.Range("N246").Value = 9500
$ws.Range("O246").Value = 10000
$ws.Range("P246").Value = 9750
$ws.Range("S246").Value = 1393

# Row 247 - Segunda
$ws.Range("M247").Value = 400
$ws.Range("N247").Value = 7500
$ws.Range("O247").Value = 8000
$ws.Range("P247").Value = 7750
$ws.Range("S247").Value = 1107
